# Convert the boolean FALSE values in H2:H27 (roboticRNAPrep column) into
# the literal text string "False", matching the author's "changing FALSE
# to False" edit.
#
# The cells currently hold the boolean value FALSE displayed via a custom
# number format ("TRUE";"TRUE";"FALSE"). We need them to become plain text
# cells (shared string "False") formatted with the built-in Text number
# format (49) using the Arial 10 font (same font already used elsewhere in
# the sheet, e.g. column I/J/...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("H2:H27")

# 1) Force the cell contents to be stored as literal text "False" rather
#    than being auto-recognised as the boolean FALSE. A leading apostrophe
#    is the standard way to tell Excel "treat this entry as text".
$target.Value = "'False"

# 2) Re-apply the formatting that's used elsewhere in the sheet (Arial 10,
#    black) by copying the format from an already-styled cell (I1 uses
#    that exact font), then switch the number format to Text (49) so the
#    string displays as typed instead of being re-interpreted.
$ws.Range("I1").Copy()
$target.PasteSpecial(-4122)
$target.NumberFormat = "@"

# 3) Update the sheet's selection to match the edited range, as recorded
#    in the saved workbook (active cell H2, selection H2:H27).
$target.Select()
